$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A57").Value = "HaUI"
$ws.Range("B57").Value = 185.8060249923673
$ws.Range("C57").Value = "2025-05-21 01:18:50"
$ws.Range("D57").Value = "/home/anodi108/Desktop/project/Do_An_Tot_Nghiep/DATN_PhamDangDong/DATN_PhamDangDong/resource/data/data_result/image_20250521_011847.jpg"
